# Apply updated placement results to the "Placings" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Placings")

# Row 14
$ws.Range("G14").Value = 415948
$ws.Range("H14").Value = "Yoshi"

# Row 15
$ws.Range("G15").Value = 592052
$ws.Range("H15").Value = "MNG Mike"

# Row 18
$ws.Range("G18").Value = 281841
$ws.Range("H18").Value = "rickbb"

# Row 19
$ws.Range("G19").Value = 592481
$ws.Range("H19").Value = "Riggs"

# Row 20
$ws.Range("F20").Value = "-"
$ws.Range("G20").Value = 340348
$ws.Range("H20").Value = "Tury"

# Row 22
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 225416
$ws.Range("H22").Value = "Wobble2King"

# Row 23
$ws.Range("D23").Value = 9
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 246986
$ws.Range("H23").Value = "Werito"

# Row 24
$ws.Range("D24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = 179962
$ws.Range("H24").Value = "Fran"

# Row 25
$ws.Range("G25").Value = 587990
$ws.Range("H25").Value = "TEC"

# Row 32
$ws.Range("D32").Value = "-"
$ws.Range("E32").Value = "-"
$ws.Range("F32").Value = "-"
$ws.Range("G32").Value = 756906
$ws.Range("H32").Value = "Rojo"

# Row 33
$ws.Range("D33").Value = 13
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 17
$ws.Range("G33").Value = 466335
$ws.Range("H33").Value = "ViviS"

# Row 35
$ws.Range("G35").Value = 1030049
$ws.Range("H35").Value = "Carreto"

# Row 36
$ws.Range("E36").Value = "-"
$ws.Range("F36").Value = 13
$ws.Range("G36").Value = 749414
$ws.Range("H36").Value = "Saru"

# Row 37
$ws.Range("E37").Value = 11
$ws.Range("F37").Value = 17
$ws.Range("G37").Value = 466863
$ws.Range("H37").Value = "Fabinni"

# Row 38
$ws.Range("G38").Value = 144909
$ws.Range("H38").Value = "Navson"

# Row 39
$ws.Range("F39").Value = 17
$ws.Range("G39").Value = 126392
$ws.Range("H39").Value = "Helsxan"

# Row 40
$ws.Range("F40").Value = "-"
$ws.Range("G40").Value = 62728
$ws.Range("H40").Value = "Marcopolo"

# Row 41
$ws.Range("G41").Value = 1030453
$ws.Range("H41").Value = "DOOM"
